$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 13 (shifts old rows 13-23 down to 14-24) ---
$ws.Rows(13).Insert()
# The freshly inserted row picks up a stray formatted-but-empty cell at A13;
# clear it so the row matches the target (no A13 cell at all).
$ws.Cells.Item(13, 1).Clear()

# --- Row 10: "Objetivos:" body text changed ---
$ws.Range("B10").Value = "Complementar os conhecimentos na Área de Operações Unitárias da Indústria Química, com aplicações na operação, análise e projeto de equipamentos."
$ws.Range("C10").Value = "Complementar os conhecimentos na Área de Operações Unitárias da Indústria Química, com aplicações na operação, análise e projeto de equipamentos."

# --- Row 13 (new): "Docentes responsáveis:" value, no label in column A ---
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("B13").Value = "8151869 - Livia Chaguri e Carvalho"
$ws.Range("C13").Value = "8151869 - Livia Chaguri e Carvalho"

# --- Row 14: "Programa resumido:" body text changed (was "Semestral") ---
$programaResumido = @"
1. Tópicos especiais de operações unitárias envolvendo fluidos. 
2. Tópicos especiais de operações unitárias envolvendo transmissão de calor e massa.
"@
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido

# --- Row 15: "Short syllabus:" no B/C values (previously wrongly carried a date) ---
$ws.Range("B15").ClearContents()
$ws.Range("C15").ClearContents()

# --- Row 16: "Programa:" body text changed to match the new syllabus text ---
$ws.Range("B16").Value = $programaResumido
$ws.Range("C16").Value = $programaResumido

# --- Row 18: "Avaliação:" no B/C values (previously wrongly carried a name) ---
$ws.Range("B18").ClearContents()
$ws.Range("C18").ClearContents()

# --- Row 19: "Método:" body text corrected ---
$metodo = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# --- Row 20: "Critério:" body text corrected ---
$ws.Range("B20").Value = "Provas e trabalhos."
$ws.Range("C20").Value = "Provas e trabalhos."

# --- Row 21: "Norma de recuperação:" body text corrected ---
$ws.Range("B21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."
$ws.Range("C21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."

# --- Row 22: "Bibliografia:" body text replaced with the full reference list ---
$bibliografia = @"
FOUST, Alan S. Princípios das Operações Unitárias. Rio de Janeiro : Guanabara Dois/LTC, 1982.
MCCABE, Warren; SMITH, Julian; HARRIOTT, Peter. Unit Operations of Chemical Engineering. Boston : McGraw-Hill, 2005.
GEANKOPLIS, Christie John. Transport Processes and Separation Process Principles. Upper Saddle River, NJ : Prentice Hall Professional Technical Reference, 2003.
COUPER, James R.; PENNEY, W. Roy; FAIR, James R.; WALAS, Stanley M. Chemical Process Equipment: Selection and Design. Amsterdam : Elsevier, c2005Boston.
PERRY, Robert H; GREEN, Don W; MALONEY, James O. Perry's Chemical Engineers' Handbook. 7th. ed. New York : McGraw-Hill, 1999.
Textos fornecidos pelo professor da disciplina
Artigos extraídos de revistas especializadas de Engenharia Química.
"@
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
